$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300290942192078
$ws.Range("B1").Value = 1.820895433425903
$ws.Range("C1").Value = 3.274490833282471
$ws.Range("D1").Value = 3.775489091873169
$ws.Range("E1").Value = 1.186664819717407
